# Reorder the header values in row 1 of the active worksheet.
# Before: A1=bedrooms_1, B1=kitchens_1, C1=living_rooms_1, D1=bedrooms_2, E1=kitchens_2, F1=living_rooms_2
# After:  A1=living_rooms_1, B1=bedrooms_1, C1=kitchens_1, D1=kitchens_2, E1=living_rooms_2, F1=bedrooms_2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "living_rooms_1"
$ws.Range("B1").Value = "bedrooms_1"
$ws.Range("C1").Value = "kitchens_1"
$ws.Range("D1").Value = "kitchens_2"
$ws.Range("E1").Value = "living_rooms_2"
$ws.Range("F1").Value = "bedrooms_2"
